$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price + 1h volume %) for rows 2-51,
# mirroring a refreshed GitHub Actions data pull.
#
# The D-column (Price) cells hold plain numeric-looking text such as
# "0.9590" or "39.41". Excel COM auto-parses a numeric-looking string
# assigned via .Value into a real Number, which would silently drop
# meaningful trailing zeros (e.g. "0.9590" -> 0.959). Force those
# specific cells to Text format first so the literal string survives,
# matching the original inlineStr text cells in the workbook. The
# E-column (Volume) strings such as "  +1.44%  " already contain
# spaces/percent signs and are never mistaken for numbers, so no
# special handling is required there.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '20.556.87'
$ws.Range('E2').Value = '  +1.44%  '
$ws.Range('D3').Value = '1.471.94'
$ws.Range('E3').Value = '  +2.01%  '
$ws.Range('E4').Value = '  +0.38%  '
$ws.Range('D5').Value = '0.9590'
$ws.Range('E5').Value = '  +4.67%  '
$ws.Range('D6').Value = '277.05'
$ws.Range('E6').Value = '  +0.83%  '
$ws.Range('D7').Value = '0.3566'
$ws.Range('E7').Value = '  -1.87%  '
$ws.Range('D8').Value = '0.3067'
$ws.Range('E8').Value = '  -0.40%  '
$ws.Range('E9').Value = '  +6.97%  '
$ws.Range('D10').Value = '39.41'
$ws.Range('E10').Value = '  +0.92%  '
$ws.Range('D11').Value = '0.06640'
$ws.Range('E11').Value = '  +2.26%  '
$ws.Range('E12').Value = '  +0.51%  '
$ws.Range('D13').Value = '5.452'
$ws.Range('E13').Value = '  +2.16%  '
$ws.Range('D14').Value = '18.08'
$ws.Range('D15').Value = '6.174'
$ws.Range('E15').Value = '  +2.27%  '
$ws.Range('D16').Value = '0.9588'
$ws.Range('E16').Value = '  +2.93%  '
$ws.Range('D17').Value = '0.00001020'
$ws.Range('E17').Value = '  +1.06%  '
$ws.Range('D18').Value = '1.474.42'
$ws.Range('E18').Value = '  +2.26%  '
$ws.Range('D19').Value = '0.05957'
$ws.Range('E19').Value = '  +6.00%  '
$ws.Range('D20').Value = '69.07'
$ws.Range('E20').Value = '  +2.52%  '
$ws.Range('D21').Value = '5.489'
$ws.Range('E21').Value = '  +1.69%  '
$ws.Range('D22').Value = '14.51'
$ws.Range('E22').Value = '  +2.26%  '
$ws.Range('D23').Value = '11.27'
$ws.Range('E23').Value = '  +4.37%  '
$ws.Range('D24').Value = '2.278'
$ws.Range('E24').Value = '  +1.99%  '
$ws.Range('D25').Value = '20.564.17'
$ws.Range('E25').Value = '  +1.34%  '
$ws.Range('D26').Value = '145.07'
$ws.Range('E26').Value = '  +5.47%  '
$ws.Range('D27').Value = '2.084'
$ws.Range('E27').Value = '  +1.82%  '
$ws.Range('D28').Value = '17.15'
$ws.Range('E28').Value = '  +1.45%  '
$ws.Range('D29').Value = '1.631.57'
$ws.Range('E29').Value = '  +2.40%  '
$ws.Range('D30').Value = '114.09'
$ws.Range('E30').Value = '  +3.63%  '
$ws.Range('D31').Value = '3.861'
$ws.Range('E31').Value = '  -2.85%  '
$ws.Range('D32').Value = '4.923'
$ws.Range('E32').Value = '  +1.93%  '
$ws.Range('D33').Value = '0.07930'
$ws.Range('E33').Value = '  +3.49%  '
$ws.Range('D34').Value = '0.7942'
$ws.Range('E34').Value = '  -0.22%  '
$ws.Range('D35').Value = '1.238'
$ws.Range('E35').Value = '  +8.97%  '
$ws.Range('D36').Value = '1.446'
$ws.Range('E36').Value = '  -1.13%  '
$ws.Range('D37').Value = '0.05740'
$ws.Range('E37').Value = '  -0.60%  '
$ws.Range('D38').Value = '4.700'
$ws.Range('E38').Value = '  +0.97%  '
$ws.Range('D39').Value = '0.02034'
$ws.Range('E39').Value = '  +2.54%  '
$ws.Range('D40').Value = '0.9594'
$ws.Range('E40').Value = '  +3.40%  '
$ws.Range('D41').Value = '10.30'
$ws.Range('E41').Value = '  +1.70%  '
$ws.Range('D42').Value = '0.1854'
$ws.Range('E42').Value = '  +0.29%  '
$ws.Range('D43').Value = '7.269'
$ws.Range('E43').Value = '  +3.65%  '
$ws.Range('D44').Value = '0.5244'
$ws.Range('E44').Value = '  +1.04%  '
$ws.Range('D45').Value = '3.508'
$ws.Range('E45').Value = '  +0.82%  '
$ws.Range('D46').Value = '12.03'
$ws.Range('E46').Value = '  +2.14%  '
$ws.Range('D47').Value = '118.57'
$ws.Range('E47').Value = '  +2.12%  '
$ws.Range('D48').Value = '0.5177'
$ws.Range('E48').Value = '  +1.38%  '
$ws.Range('D49').Value = '1.800'
$ws.Range('E49').Value = '  +4.32%  '
$ws.Range('D50').Value = '0.06434'
$ws.Range('E50').Value = '  +0.54%  '
$ws.Range('D51').Value = '0.9906'
$ws.Range('E51').Value = '  +1.43%  '
